# Actualización desde MV -datos-
# Append the next 6 daily rows (09-10-2021 .. 14-10-2021) to Sheet1,
# replicating the values of the last existing row (12836 / 266 / 393).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @("09-10-2021", "10-10-2021", "11-10-2021", "12-10-2021", "13-10-2021", "14-10-2021")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1
$startRow = $lastRow + 1

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $cellA = $ws.Cells.Item($row, 1)

    # Column A stores the dates as plain text (e.g. "08-10-2021"), not
    # real Excel date serials. Assigning the literal string straight to
    # .Value lets Excel's date auto-recognition kick in and silently
    # turn it into a date number. Routing it through TEXT(...,"@") keeps
    # it a genuine string, then Copy / PasteSpecial values-only drops the
    # formula and leaves a plain text cell with no added number format,
    # matching the rest of the sheet (which carries no explicit style on
    # its data cells).
    $cellA.Formula = '=TEXT("' + $dates[$i] + '","@")'
    $cellA.Copy()
    $cellA.PasteSpecial(-4163)

    $ws.Cells.Item($row, 2).Value = 12836
    $ws.Cells.Item($row, 3).Value = 266
    $ws.Cells.Item($row, 4).Value = 393
}
$excel.CutCopyMode = $false
